$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
$pkgFooter = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# Edit 1: "North Central University" -> split into "North" / "c" / "entral University"
#         (net text becomes "Northcentral University" - corrects the school name spelling)
# ---------------------------------------------------------------------
$para1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "North Central University") {
        $para1 = $p
        break
    }
}
if ($null -eq $para1) {
    throw "Could not find the 'North Central University' paragraph"
}
# Range spanning the whole paragraph's visible text (excludes the trailing paragraph mark)
$start1 = $para1.Range.Start
$end1 = $para1.Range.End - 1

$target1 = $d.Range($start1, $end1)
$xml1 = $pkgHeader + '<w:r><w:t>North</w:t></w:r><w:r><w:t>c</w:t></w:r><w:r><w:t>entral University</w:t></w:r>' + $pkgFooter
$target1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Edit 2: "driving ability the dependent result" -> "driving ability, the dependent result"
#         (insert a comma, splitting the run in three: before-comma / comma / after-comma)
#         The paragraph also has an unchanged trailing run ("  ") that must be preserved
#         in place, so the whole remainder of the paragraph is rebuilt together.
# ---------------------------------------------------------------------
$para2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*driving ability the dependent*") {
        $para2 = $p
        break
    }
}
if ($null -eq $para2) {
    throw "Could not find the paragraph containing 'driving ability the dependent'"
}

$find2 = $d.Content.Find
$find2.ClearFormatting()
$found2 = $find2.Execute(" the specific problem of research.  Harris (2008) uses an example of investigating the " + [char]8220 + "influence of music on driving ability." + [char]8221 + "  In this situation, music is the independent variable versus driving ability the dependent result.  However, music could refer to either the categorical genre or perhaps the volume level.  Meanwhile, the driver" + [char]8217 + "s ability could be a measurement of maintaining lane alignment or parallel parking.  These nuances to the question have a substantial impact on all aspects of the design and must be declared upfront.",
                          $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the target run text for edit 2"
}
$start2 = $find2.Parent.Start
# extend to the true end of the paragraph (excludes the trailing paragraph mark) so any
# unchanged trailing runs (e.g. the "  " run) stay put instead of being shuffled by InsertXML
$end2 = $para2.Range.End - 1

$runA = " the specific problem of research.  Harris (2008) uses an example of investigating the " + [char]8220 + "influence of music on driving ability." + [char]8221 + "  In this situation, music is the independent variable versus driving ability"
$runC = " the dependent result.  However, music could refer to either the categorical genre or perhaps the volume level.  Meanwhile, the driver" + [char]8217 + "s ability could be a measurement of maintaining lane alignment or parallel parking.  These nuances to the question have a substantial impact on all aspects of the design and must be declared upfront."

$target2 = $d.Range($start2, $end2)
$xml2 = $pkgHeader + '<w:r><w:t xml:space="preserve">' + $runA + '</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve">' + $runC + '</w:t></w:r><w:r w:rsidR="00161929"><w:t xml:space="preserve">  </w:t></w:r>' + $pkgFooter
$target2.InsertXML($xml2)
